# % AS coverage code and it's results are added
# Adds a new data row (17) to Sheet1 with the "ALL (458)" / "5888(with% Ases
# covered part)" pair, matching the formatting already used for the
# analogous cells (D column label style, and the B-column numeric/result
# style), and extends the line chart's category/value series ranges to
# include the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- copy formatting for the new cells from equivalent existing cells ---
# A17 should look like the other "label" cells in column D (style s="1").
$ws.Range("D16").Copy()
$ws.Range("A17").PasteSpecial(-4122)

# B17 should look like the later numeric/result cells in column B (style s="2").
$ws.Range("B6").Copy()
$ws.Range("B17").PasteSpecial(-4122)

# --- set the new cell values ---
$ws.Range("A17").Value = "ALL (458)"
$ws.Range("B17").Value = "5888(with% Ases covered part)"

# --- extend the chart's series ranges to include the new row ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.SeriesCollection(1).Formula = "=SERIES(Sheet1!`$B`$1,Sheet1!`$A`$2:`$A`$17,Sheet1!`$B`$2:`$B`$17,1)"
